# The crawler re-scraped the "electronics" listing later the same day.
# Two things changed as a result:
#   1) the "timestamp" column (O) moves from 07:05:09 to 21:00:12 for every row.
#   2) a few of the battery products (originally in rows 61-63/65-66) came back
#      from the site in a different order / with refreshed price & rating data;
#      row 64 (Varta Longlife Max Power C) was unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update timestamp column O for every data row (2-86)
$ws.Range("O2:O86").Value = "2022-08-01 21:00:12"

# 2) Rewrite rows 61-63 and 65-66 with the re-crawled product data.
#    Text cells are forced to the "@" (text) format before assignment so that
#    numeric-looking strings (ids, prices, ...) are stored as text, matching
#    the rest of the sheet; the temporary format is cleared again afterwards
#    so no stray cell styles are left behind. Columns E/F (ratingAmount /
#    ratingValue) are genuinely numeric and are set as plain numbers.

# Row 61
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = "4119046"
$ws.Range("A61").ClearFormats()
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = "Varta Ultra Lithium AA 4er Bli"
$ws.Range("B61").ClearFormats()
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/aa/varta-ultra-lithium-aa-4er-bli/p/4119046"
$ws.Range("C61").ClearFormats()
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = "4ST"
$ws.Range("D61").ClearFormats()
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 5
$ws.Range("G61").NumberFormat = "@"
$ws.Range("G61").Value = "Varta"
$ws.Range("G61").ClearFormats()
$ws.Range("H61").NumberFormat = "@"
$ws.Range("H61").Value = "14.95"
$ws.Range("H61").ClearFormats()
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "3.74/1ST"
$ws.Range("I61").ClearFormats()
$ws.Range("J61").NumberFormat = "@"
$ws.Range("J61").Value = "Preis pro 1 Stück"
$ws.Range("J61").ClearFormats()
$ws.Range("K61").NumberFormat = "@"
$ws.Range("K61").Value = "3.74"
$ws.Range("K61").ClearFormats()
$ws.Range("L61").NumberFormat = "@"
$ws.Range("L61").Value = "1ST"
$ws.Range("L61").ClearFormats()
$ws.Range("M61").NumberFormat = "@"
$ws.Range("M61").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'aa']"
$ws.Range("M61").ClearFormats()
$ws.Range("N61").NumberFormat = "@"
$ws.Range("N61").Value = "Varta Ultra Lithium AA 4er Bli 14.95 Schweizer Franken"
$ws.Range("N61").ClearFormats()

# Row 62
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "4905486"
$ws.Range("A62").ClearFormats()
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = "Alkaline Batterie 3LR12/4.5V"
$ws.Range("B62").ClearFormats()
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/alkaline-batterie-3lr1245v/p/4905486"
$ws.Range("C62").ClearFormats()
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "1ST"
$ws.Range("D62").ClearFormats()
$ws.Range("E62").Value = 1
$ws.Range("F62").Value = 1
$ws.Range("G62").NumberFormat = "@"
$ws.Range("G62").Value = "Coop"
$ws.Range("G62").ClearFormats()
$ws.Range("H62").NumberFormat = "@"
$ws.Range("H62").Value = "5.95"
$ws.Range("H62").ClearFormats()
$ws.Range("I62").NumberFormat = "@"
$ws.Range("I62").Value = "5.95/1ST"
$ws.Range("I62").ClearFormats()
$ws.Range("J62").NumberFormat = "@"
$ws.Range("J62").Value = "Preis pro 1 Stück"
$ws.Range("J62").ClearFormats()
$ws.Range("K62").NumberFormat = "@"
$ws.Range("K62").Value = "5.95"
$ws.Range("K62").ClearFormats()
$ws.Range("L62").NumberFormat = "@"
$ws.Range("L62").Value = "1ST"
$ws.Range("L62").ClearFormats()
$ws.Range("M62").NumberFormat = "@"
$ws.Range("M62").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("M62").ClearFormats()
$ws.Range("N62").NumberFormat = "@"
$ws.Range("N62").Value = "Alkaline Batterie 3LR12/4.5V 5.95 Schweizer Franken"
$ws.Range("N62").ClearFormats()

# Row 63
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "6986541"
$ws.Range("A63").ClearFormats()
$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = "Duracell Batterie (CR2032, 4 Stück)"
$ws.Range("B63").ClearFormats()
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/duracell-batterie-cr2032-4-stueck/p/6986541"
$ws.Range("C63").ClearFormats()
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "4ST"
$ws.Range("D63").ClearFormats()
$ws.Range("E63").Value = ""
$ws.Range("F63").Value = 0
$ws.Range("G63").NumberFormat = "@"
$ws.Range("G63").Value = "Duracell"
$ws.Range("G63").ClearFormats()
$ws.Range("H63").NumberFormat = "@"
$ws.Range("H63").Value = "8.95"
$ws.Range("H63").ClearFormats()
$ws.Range("I63").NumberFormat = "@"
$ws.Range("I63").Value = "2.24/1ST"
$ws.Range("I63").ClearFormats()
$ws.Range("J63").NumberFormat = "@"
$ws.Range("J63").Value = "Preis pro 1 Stück"
$ws.Range("J63").ClearFormats()
$ws.Range("K63").NumberFormat = "@"
$ws.Range("K63").Value = "2.24"
$ws.Range("K63").ClearFormats()
$ws.Range("L63").NumberFormat = "@"
$ws.Range("L63").Value = "1ST"
$ws.Range("L63").ClearFormats()
$ws.Range("M63").NumberFormat = "@"
$ws.Range("M63").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("M63").ClearFormats()
$ws.Range("N63").NumberFormat = "@"
$ws.Range("N63").Value = "Duracell Batterie (CR2032, 4 Stück) 43% Aktion 8.95 Schweizer Franken statt 15.90 Schweizer Franken"
$ws.Range("N63").ClearFormats()

# Row 65
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "4096628"
$ws.Range("A65").ClearFormats()
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = "Varta Longlife Power Batterien AAA/LR03 6 Stück"
$ws.Range("B65").ClearFormats()
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "/de/haushalt-tier/elektroartikel-batterien/batterien/andere-batterien/varta-longlife-power-batterien-aaalr03-6-stueck/p/4096628"
$ws.Range("C65").ClearFormats()
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "6ST"
$ws.Range("D65").ClearFormats()
$ws.Range("E65").Value = ""
$ws.Range("F65").Value = 0
$ws.Range("G65").NumberFormat = "@"
$ws.Range("G65").Value = "Varta"
$ws.Range("G65").ClearFormats()
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = "12.95"
$ws.Range("H65").ClearFormats()
$ws.Range("I65").NumberFormat = "@"
$ws.Range("I65").Value = "2.16/1ST"
$ws.Range("I65").ClearFormats()
$ws.Range("J65").NumberFormat = "@"
$ws.Range("J65").Value = "Preis pro 1 Stück"
$ws.Range("J65").ClearFormats()
$ws.Range("K65").NumberFormat = "@"
$ws.Range("K65").Value = "2.16"
$ws.Range("K65").ClearFormats()
$ws.Range("L65").NumberFormat = "@"
$ws.Range("L65").Value = "1ST"
$ws.Range("L65").ClearFormats()
$ws.Range("M65").NumberFormat = "@"
$ws.Range("M65").Value = "['haushalt-tier', 'elektroartikel-batterien', 'batterien', 'andere-batterien']"
$ws.Range("M65").ClearFormats()
$ws.Range("N65").NumberFormat = "@"
$ws.Range("N65").Value = "Varta Longlife Power Batterien AAA/LR03 6 Stück 12.95 Schweizer Franken"
$ws.Range("N65").ClearFormats()

# Row 66
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "6378875"
$ws.Range("A66").ClearFormats()
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = "Wilkinson Intuition 4in1 Finish Styler"
$ws.Range("B66").ClearFormats()
$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "/de/haushalt-tier/elektroartikel-batterien/elektrogeraete/beautygeraete/wilkinson-intuition-4in1-finish-styler/p/6378875"
$ws.Range("C66").ClearFormats()
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "1ST"
$ws.Range("D66").ClearFormats()
$ws.Range("E66").Value = ""
$ws.Range("F66").Value = 0
$ws.Range("G66").NumberFormat = "@"
$ws.Range("G66").Value = "Wilkinson"
$ws.Range("G66").ClearFormats()
$ws.Range("H66").NumberFormat = "@"
$ws.Range("H66").Value = "29.95"
$ws.Range("H66").ClearFormats()
$ws.Range("I66").NumberFormat = "@"
$ws.Range("I66").Value = "29.95/1ST"
$ws.Range("I66").ClearFormats()
$ws.Range("J66").NumberFormat = "@"
$ws.Range("J66").Value = "Preis pro 1 Stück"
$ws.Range("J66").ClearFormats()
$ws.Range("K66").NumberFormat = "@"
$ws.Range("K66").Value = "29.95"
$ws.Range("K66").ClearFormats()
$ws.Range("L66").NumberFormat = "@"
$ws.Range("L66").Value = "1ST"
$ws.Range("L66").ClearFormats()
$ws.Range("M66").NumberFormat = "@"
$ws.Range("M66").Value = "['haushalt-tier', 'elektroartikel-batterien', 'elektrogeraete', 'beautygeraete']"
$ws.Range("M66").ClearFormats()
$ws.Range("N66").NumberFormat = "@"
$ws.Range("N66").Value = "Wilkinson Intuition 4in1 Finish Styler 29.95 Schweizer Franken"
$ws.Range("N66").ClearFormats()
